# Daily attendance processing - 2026-01-13 04:31:33
# Normalize the "Recorded By" (column G) lists on the Session Analysis
# Results sheet: pull the "System" entry (exact case) to the front of the
# list when present, otherwise sort the recorder names alphabetically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = @($text -split ', ')

    if ($parts.Count -le 1) {
        continue
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $rest = @($parts | Where-Object { -not $_.Equals("System") })
        $newParts = @("System") + $rest
    } else {
        $newParts = @($parts | Sort-Object)
    }

    $newText = ($newParts -join ", ")

    if (-not $newText.Equals($text)) {
        $cell.Value = $newText
    }
}
